$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for 945353e5-...md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-15 17:22:11"

# Sheet "zh-cn": Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-15 17:21:57"
$wsZhCn.Range("K2").Value = "2016-11-15 17:22:54"

# Sheet "de-de": Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-11-15 17:23:16"
